# Apply updated "想去人数" (want-to-go count) values per the diff.
# (regenerated gh-pages data snapshot at commit 456a3b4)
# Sheet order in workbook: 1=展览, 2=演出, 3=本地生活, 4=全部类型
$wb = $excel.ActiveWorkbook

# Worksheet 1 (展览)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 2682
$ws.Range("F7").Value = 2190
$ws.Range("F8").Value = 1807
$ws.Range("F11").Value = 2463
$ws.Range("F13").Value = 229
$ws.Range("F17").Value = 108
$ws.Range("F18").Value = 9129
$ws.Range("F20").Value = 7075
$ws.Range("F21").Value = 11572
$ws.Range("F24").Value = 234
$ws.Range("F25").Value = 341
$ws.Range("F26").Value = 554
$ws.Range("F27").Value = 2553
$ws.Range("F30").Value = 2490
$ws.Range("F31").Value = 665
$ws.Range("F33").Value = 4498
$ws.Range("F34").Value = 855
$ws.Range("F35").Value = 347
$ws.Range("F36").Value = 40
$ws.Range("F37").Value = 515

# Worksheet 2 (演出)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 116
$ws.Range("F14").Value = 65

# Worksheet 3 (本地生活)
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 624

# Worksheet 4 (全部类型)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 624
$ws.Range("F5").Value = 2682
$ws.Range("F9").Value = 2190
$ws.Range("F11").Value = 1807
$ws.Range("F14").Value = 2463
$ws.Range("F15").Value = 116
$ws.Range("F17").Value = 229
$ws.Range("F21").Value = 108
$ws.Range("F22").Value = 9129
$ws.Range("F24").Value = 7075
$ws.Range("F25").Value = 11572
$ws.Range("F28").Value = 234
$ws.Range("F29").Value = 341
$ws.Range("F31").Value = 554
$ws.Range("F33").Value = 2553
$ws.Range("F39").Value = 4498
$ws.Range("F40").Value = 65
$ws.Range("F46").Value = 515

